# DoR.docx edit: replace the old "Detailed Web Magazine Requirements" outline
# with the new Definition-of-Ready (INVEST-style) criteria for the
# "Srsťoplsť" web magazine project (Jasnost, Dosáhnutelnost, Testovatelnost,
# Hodnota, Nezávislost, Vyjednatelnost).
#
# The new content mixes plain runs with <w:proofErr> spell/grammar-check
# markers around a few words/names (e.g. "Srsťoplsť"), so rather than doing a
# series of Find/Replace calls (which would not recreate the proofErr wrapper
# runs) we build the exact target run/paragraph structure as OOXML and push
# it into the document body in one shot via Range.InsertXML, which is the
# Word COM-interop method for inserting raw WordprocessingML.

$d = $word.ActiveDocument

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p><w:r><w:t>Jasnost: Každý úkol spojený s webem, ať už se týká obsahu, designu nebo funkcí, musí být jasně popsaný tak, aby bylo zřejmé, co má být dosaženo. Například, příběh o adopci zvířat by měl mít definovaný rozsah, strukturu článku a všechny potřebné zdroje.</w:t></w:r></w:p>
          <w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Dosáhnutelnost</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: Úkoly musí být rozloženy tak, aby bylo možné je dokončit v rámci jednoho sprintu. Pro "</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Srsťoplsť</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>" to může znamenat například vytvoření jedné rubriky o péči o zvířata nebo přidání nové sekce recenzí produktů pro zvířata.</w:t></w:r></w:p>
          <w:p><w:r><w:t>Testovatelnost: Pro každý úkol by měla existovat sada testovacích kritérií. Například, je-li nová funkce interaktivní mapy obchodů se zvířaty přidána na web, měly by být specifikovány testy pro její funkčnost na různých zařízeních a prohlížečích.</w:t></w:r></w:p>
          <w:p><w:r><w:t>Hodnota: Každý úkol by měl přinést zjevný přínos pro čtenáře "</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Srsťoplsť</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>", například poskytnutím užitečných informací nebo zlepšením uživatelské zkušenosti na webu.</w:t></w:r></w:p>
          <w:p><w:r><w:t>Nezávislost: Úkoly by měly být definovány tak, aby na sebe nebyly navzájem závislé a bylo možné na nich pracovat paralelně. To znamená, že například přidání nové fotogalerie nesmí být závislé na přepracování domovské stránky.</w:t></w:r></w:p>
          <w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Vyjednatelnost</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: Detaily úkolu by měly být dostatečně otevřené pro diskuzi a úpravy na základě zpětné vazby od redakčního týmu "</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Srsťoplsť</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">" a od </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>čtenářů.Tyto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> kritéria zajistí, že před zahájením práce na sprintu je každý úkol pro web magazínu "</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Srsťoplsť</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>" dobře pochopen, zpracovatelný a schválený všemi členy týmu.</w:t></w:r></w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

# Replace the whole document body (Content = body minus the final paragraph
# mark/sectPr) with the new paragraphs; page setup (sectPr) is untouched.
$d.Content.InsertXML($xml)

Write-Host "DoR.docx updated: body now has $($d.Paragraphs.Count) paragraphs"
